# "Retouche de 2-3 maquette + documentation avec maquette"
#
# Fills in the two previously-blank journal rows (19 & 20) on the only
# worksheet, adding two new shared strings along the way, and leaves the
# selection on E20 (matching the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: new task entry with date + duration, description left blank.
$ws.Range("B19").Value = "Ajout des maquettes et commentaires dans la documentation"
$ws.Range("C19").Value = "5/6/2021"
$ws.Range("D19").Value = 60

# Row 20: task name + duration + description filled in; date left blank.
$ws.Range("B20").Value = "Documentation"
$ws.Range("D20").Value = 120
$ws.Range("E20").Value = "Debut de usercase + test"

# Row 19 grew taller to fit the wrapped task description.
$ws.Rows.Item(19).RowHeight = 45

# Leave the selection where the author left it.
$ws.Range("E20").Select()
